# Updated Test Cases List: revised comments / retest notes across several
# rows and refreshed view state (scroll position + active cell selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 13: variant reads (Test Done / Name already "?", left as-is by value,
#     but rewritten here to be explicit/idempotent) ---
$ws.Range("C13").Value = "?"
$ws.Range("D13").Value = "?"

# --- Row 20: seq ont transcript ---
$ws.Range("C20").Value = "?"
$ws.Range("D20").Value = "?"
$ws.Range("E20").Value = "retest, look up in refsew, ensembl"

# --- Row 21: seq ont all transcripts ---
$ws.Range("C21").Value = "?"
$ws.Range("D21").Value = "?"
$ws.Range("E21").Value = "retest, look up in refsew, ensembl"

# --- Row 38: target ---
$ws.Range("E38").Value = "CRAVAT db not matching any available dbs. xls files in testing folder"

# --- Row 44: pubmed articles ---
$ws.Range("E44").Value = "test just for a ballpark figure"

# --- Row 45: pubmed search term ---
$ws.Range("C45").Value = "?"
$ws.Range("D45").Value = "?"
$ws.Range("E45").Value = "all vars with pubmed hits should have a pubmed link"

# --- Row 46: 1000 genomes freq ---
$ws.Range("C46").Value = "?"
$ws.Range("D46").Value = "?"
$ws.Range("E46").Value = "browser.1000genomes.com"

# --- Rows 47-56: esp6500 freq / ExAC freq * -> Name column now "pop_stats" ---
$ws.Range("D47").Value = "pop_stats"
$ws.Range("E47").Value = "out of date. also, when is it null vs when is it zero? see uid CYP19A1 vs uid CYP19A1_NC in pop_stats"
$ws.Range("D48").Value = "pop_stats"
$ws.Range("D49").Value = "pop_stats"
$ws.Range("D50").Value = "pop_stats"
$ws.Range("D51").Value = "pop_stats"
$ws.Range("D52").Value = "pop_stats"
$ws.Range("D53").Value = "pop_stats"
$ws.Range("D54").Value = "pop_stats"
$ws.Range("D55").Value = "pop_stats"
$ws.Range("D56").Value = "pop_stats"

# --- Row 57: mupit link ---
$ws.Range("C57").Value = "?"
$ws.Range("D57").Value = "?"
$ws.Range("E57").Value = "regtest for now, use genes available in mupit"

# --- Row 58: in tcga mutation cluster ---
$ws.Range("C58").Value = "?"
$ws.Range("D58").Value = "?"
$ws.Range("E58").Value = "regtest for now, use mupit to find aas in hotspots on certain genes"

# --- Row 59: identical samps in study (new rounding test case row) ---
$ws.Range("C59").Value = "?"
$ws.Range("D59").Value = "?"
$ws.Range("E59").Value = "include in vcf testing?"

# --- Refresh the saved view state: scroll the window so row 10 is at the
#     top, then select E21 (matches the author's editing position). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("E21").Select() | Out-Null
